$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.813.85'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.08%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.469.14'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.45%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.66'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.90'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +4.64%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.480'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.30%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.38%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.22%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +4.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.067.57'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.89'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +5.11%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.477.93'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000171'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.851.94'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.36'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +4.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.41'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +5.97%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.28'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +3.39%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '389.65'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.565'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '74.89'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.36%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.617.19'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.55%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.179'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -6.89%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +6.13%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.19'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.89%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.40%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.80%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.78'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.43%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.27'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +4.80%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '170.49'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.15%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +7.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '31.23'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +19.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.510.32'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.62%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0769'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.43%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.800'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.84%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.22'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.65%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +3.67%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +3.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.597.55'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +5.86%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.97%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.24'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +11.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.78'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.08%  '
$ws.Range("B51").Value = 'FirstDigitalUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.04%  '
